$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 14:52"

# 2. Update Alemania (row 8) statistics
$ws.Range("B8").Value = 158768
$ws.Range("C8").Value = 10
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 6136

# 3. Arabia Saudita overtakes Irlanda and Suecia in the ranking (rows 23-25)
#    Row 23 -> Arabia Saudita, Row 24 -> Irlanda, Row 25 -> Suecia
$ws.Range("A23").Value = "Arabia Saudita"
$ws.Range("B23").Value = 20077
$ws.Range("C23").Value = 1266
$ws.Range("D23").Value = 2784
$ws.Range("E23").Value = 17141
$ws.Range("F23").Value = 118
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 152

$ws.Range("A24").Value = "Irlanda"
$ws.Range("B24").Value = 19648
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 9233
$ws.Range("E24").Value = 9313
$ws.Range("F24").Value = 146
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 1102

$ws.Range("A25").Value = "Suecia"
$ws.Range("B25").Value = 19621
$ws.Range("C25").Value = 695
$ws.Range("D25").Value = 1005
$ws.Range("E25").Value = 16261
$ws.Range("F25").Value = 524
$ws.Range("G25").Value = 81
$ws.Range("H25").Value = 2355

# 4. Kuwait (row 60): fix swapped "Casos activos" / "Recuperados" values
$ws.Range("D60").Value = 1176
$ws.Range("E60").Value = 2241

# 5. Kenia overtakes Jamaica in the ranking (rows 117-118)
#    Row 117 -> Kenia, Row 118 -> Jamaica
$ws.Range("A117").Value = "Kenia"
$ws.Range("B117").Value = 374
$ws.Range("C117").Value = 11
$ws.Range("D117").Value = 124
$ws.Range("E117").Value = 236
$ws.Range("F117").Value = 2
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 14

$ws.Range("A118").Value = "Jamaica"
$ws.Range("B118").Value = 364
$ws.Range("C118").Value = 14
$ws.Range("D118").Value = 29
$ws.Range("E118").Value = 328
$ws.Range("F118").Value = 3
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 7
